# batch-template.xlsx: add new export columns (Contributor, Language, Rights,
# Owner, ExpirationDate, Tags:ROBOTS, Tags:publishing_entity) to the existing
# test-file metadata table, reordering so the original "-Path" formula moves
# to column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the whole table layout is changing (columns are being
# inserted and re-ordered), so wipe old cell content/formatting first.
$ws.Cells.Clear()

# ---- column widths (engine stores width_xml = ColumnWidth + 5/6, quantized
# to 1/6 char) -------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 15.307291666666666   # -> 16.140625
$ws.Columns("B").ColumnWidth = 30.592447916666668   # -> 31.42578125
$ws.Columns("C").ColumnWidth = 30.736979166666668   # -> 31.5703125
$ws.Columns("D").ColumnWidth = 11.736979166666666   # -> 12.5703125
$ws.Columns("E").ColumnWidth = 20.592447916666668   # -> 21.42578125
$ws.Columns("F").ColumnWidth = 13.451822916666666   # -> 14.28515625
$ws.Columns("G").ColumnWidth = 10.451822916666666   # -> 11.28515625
$ws.Columns("H").ColumnWidth = 8.451822916666666    # -> 9.28515625
$ws.Columns("I").ColumnWidth = 20.022135416666668   # -> 20.85546875
$ws.Columns("J").ColumnWidth = 22.877604166666668   # -> 23.7109375
$ws.Columns("K").ColumnWidth = 14.451822916666666   # -> 15.28515625

# ---- original A/B/C/F columns + row layout (these strings already exist
# in the shared-string table, so fill order here doesn't affect new indices)
$ws.Range("A1").Formula = '=T("-Path")'
$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "Description"
$ws.Range("F1").Value = "Creator"

$ws.Range("A2").Value = "test.pdf"
$ws.Range("B2").Value = "Gwent: The Witcher Card Game"
$ws.Range("C2").Value = "This is a test to see if I can batch edit metadata"
$ws.Range("F2").Value = "Von Haynes"
$ws.Rows(2).RowHeight = 30

$ws.Range("A3").Value = "nested\test.pdf"
$ws.Range("B3").Value = "Destiny 2 Coming to PC!"
$ws.Range("C3").Value = "This is a test to see if I can convert nested files"
$ws.Range("F3").Value = "Jumoke Hodari"
$ws.Rows(3).RowHeight = 30

$ws.Range("A4").Value = "test.docx"
$ws.Range("B4").Value = "Hearthstone Tournament"
$ws.Range("C4").Value = "This is a test to see if I can modify word documents"
$ws.Range("F4").Value = "James Haynes"
$ws.Rows(4).RowHeight = 30

$ws.Range("A5").Value = "test.xlsx"
$ws.Range("B5").Value = "Elder Scrolls Legends now on Ipad"
$ws.Range("C5").Value = "This is a test to see if I can modify excel docs"
$ws.Range("F5").Value = "Jane Doe"
$ws.Rows(5).RowHeight = 30

$ws.Range("A6").Value = "nested\test.pptx"
$ws.Range("B6").Value = "Now 10% more awesome"
$ws.Range("C6").Value = "Testing for nested pptx"
$ws.Range("F6").Value = "Van Jones"

# ---- new columns G:K (Contributor/Language/Rights/Owner/ExpirationDate)
# filled top-to-bottom, left-to-right so the shared-string table is built in
# the same order the original authoring tool produced it in -----------
$ws.Range("G1").Value = "Contributor"
$ws.Range("H1").Value = "Language"
$ws.Range("I1").Value = "Rights"
$ws.Range("J1").Value = "Owner"
$ws.Range("K1").Value = "ExpirationDate"

$ws.Range("H2").Value = "en"
$ws.Range("I2").Value = "Copyright &copy; 2017"
$ws.Range("J2").Value = "American Bar Association"
$ws.Range("K2").Value = 42874.135416666664

$ws.Range("H3").Value = "en"
$ws.Range("I3").Value = "Copyright &copy; 2018"
$ws.Range("J3").Value = "American Bar Association"
$ws.Range("K3").Value = 42875.763194444444

$ws.Range("H4").Value = "en"
$ws.Range("I4").Value = "Copyright &copy; 2019"
$ws.Range("J4").Value = "American Bar Association"
$ws.Range("K4").Value = 42874.135416666664

$ws.Range("H5").Value = "en"
$ws.Range("I5").Value = "Copyright &copy; 2020"
$ws.Range("J5").Value = "American Bar Association"
$ws.Range("K5").Value = 42875.135416608799

$ws.Range("H6").Value = "en"
$ws.Range("I6").Value = "Copyright &copy; 2021"
$ws.Range("J6").Value = "American Bar Association"
$ws.Range("K6").Value = 42876.135416608799

# ---- new column D (Tags:ROBOTS), then E (Tags:publishing_entity) --------
$ws.Range("D1").Value = "Tags:ROBOTS"
$ws.Range("D2").Value = "FOLLOW"
$ws.Range("D3").Value = "FOLLOW"
$ws.Range("D4").Value = "FOLLOW"
$ws.Range("D5").Value = "FOLLOW"
$ws.Range("D6").Value = "FOLLOW"

$ws.Range("E1").Value = "Tags:publishing_entity"
$ws.Range("E2").Value = "PT"
$ws.Range("E3").Value = "PT"
$ws.Range("E4").Value = "PT"
$ws.Range("E5").Value = "PT"
$ws.Range("E6").Value = "PT"

# ---- styles: Title column gets an (empty) alignment xf, Description wraps,
# ExpirationDate gets the custom date number format -----------------------
$ws.Range("C1:C6").WrapText = $true
$ws.Range("B1:B6").HorizontalAlignment = 1
$ws.Range("K1:K6").NumberFormat = 'yyyy:mm:dd\ hh:mm'

# ---- selection -------------------------------------------------------
$ws.Range("E11").Select()
